# Update the cryptocurrency price ("D") and 1h volume change ("E") columns
# for the rows whose figures changed in this refresh of the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "42.091.70"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.240.76"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.16"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.25"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.15"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0954"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.94"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "2.576.75"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.37"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.839"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "2.239.16"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "42.018.02"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.70"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.23"
$ws.Range("E22").Value = "  +8.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.12"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.41"
$ws.Range("E26").Value = "  -3.29%  "
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.70"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.60"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.61"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0802"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.93"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.124"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("E36").Value = "  -6.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.32"
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("E38").Value = "  -1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.17"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.83"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("E43").Value = "  -1.16%  "
$ws.Range("E44").Value = "  -2.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.92"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "2.450.74"
$ws.Range("E51").Value = "  +0.15%  "
